# "Generate Report for Handback" — mark the localization status sheets as
# handed back: update the Status text, fill in the Latest Target File
# (hyperlinked) / Latest Handback File / Latest Handback DateTime columns
# for zh-cn and de-de, and widen the columns that now hold the longer
# filenames and status text.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$mdFileName       = "743b15eb-139e-47e3-840b-579365822d91.md"
$mdHyperlinkAddr  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5939d601455e826eae4063ee8c0afe16e61ba3d6/e2e/743b15eb-139e-47e3-840b-579365822d91.md"

# Column widths in this engine are quantized to 1/6 of a character when
# persisted, so the literal inputs below are the values that round-trip to
# the intended stored widths (~30 and 40 character units respectively).
$widthWide = 29.166666666666668   # -> stored width 30 (closest reachable to 29.9777047293527)
$widthForty = 39.166666666666664  # -> stored width 40

# ---------------------------------------------------------------------
# Overview sheet: the "zh-cn" / "de-de" status columns (E2/F2) flip from
# "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $widthWide
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $widthWide

# ---------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de both receive the same kind of
# update, just with locale-specific handback file names / timestamps.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdHyperlinkAddr, "", "", $mdFileName)
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "743b15eb-139e-47e3-840b-579365822d91.b2d28a441e03d7be4e2c767a5fb99be0ce9734a1.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-18 11:02:17"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $widthWide
$wsZhCn.Range("I1").EntireColumn.ColumnWidth = $widthForty
$wsZhCn.Range("J1").EntireColumn.ColumnWidth = $widthForty

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdHyperlinkAddr, "", "", $mdFileName)
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "743b15eb-139e-47e3-840b-579365822d91.b2d28a441e03d7be4e2c767a5fb99be0ce9734a1.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-18 11:02:25"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $widthWide
$wsDeDe.Range("I1").EntireColumn.ColumnWidth = $widthForty
$wsDeDe.Range("J1").EntireColumn.ColumnWidth = $widthForty

Write-Host "Handback report generated."
